# Add a new log row (row 13) to the maintenance log sheet, mirroring the
# layout of the existing rows. Columns A/E/H/J/K hold digit-only values that
# must stay TEXT (not be auto-coerced to numbers) while keeping the cells on
# the sheet's default style - so for those we briefly switch the cell to the
# "Text" number format, assign the value, then restore the "Normal" style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 13

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Cells.Item($row, 1) "1306"

$ws.Cells.Item($row, 2).Value = 45540
$ws.Cells.Item($row, 2).NumberFormat = "yyyy-mm-dd"

$ws.Cells.Item($row, 3).Value = 45540
$ws.Cells.Item($row, 3).NumberFormat = "yyyy-mm-dd"

$ws.Cells.Item($row, 4).Value = "cesar ramirez"

Set-TextValue $ws.Cells.Item($row, 5) "5587964476"

$ws.Cells.Item($row, 6).Value = "UCL"
$ws.Cells.Item($row, 7).Value = "Refrigeracion Liquida"

Set-TextValue $ws.Cells.Item($row, 8) "14123"

$ws.Cells.Item($row, 9).Value = "NZXT"

Set-TextValue $ws.Cells.Item($row, 10) "123412"
Set-TextValue $ws.Cells.Item($row, 11) "54324"

$ws.Cells.Item($row, 12).Value = "No"
$ws.Cells.Item($row, 13).Value = "No"
$ws.Cells.Item($row, 14).Value = "No"
$ws.Cells.Item($row, 15).Value = "No"
$ws.Cells.Item($row, 16).Value = "Sí"
$ws.Cells.Item($row, 17).Value = "Sí"
$ws.Cells.Item($row, 18).Value = "El equipo presenta falla por oxidacion y fuga de liquido"
$ws.Cells.Item($row, 19).Value = "Correctivo, Otro"
$ws.Cells.Item($row, 20).Value = "Se procedio a reparar la carcasa y sellar los tubos del radiador"
$ws.Cells.Item($row, 21).Value = "Si"
$ws.Cells.Item($row, 22).Value = "Alcohol Isopropílico, Aislantes, Liquido Limpiador Multiusos"
$ws.Cells.Item($row, 23).Value = "Juan Daniel Ramírez Zamora"
$ws.Cells.Item($row, 24).Value = "cesar ramirez"
